$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Break apart the merged B11:B12 / C11:C12 / E11:E12 / F11:F12 cells so
#    that each row (11, 12, 13) can hold its own independent values again.
# ---------------------------------------------------------------------------
$ws.Range("B11:B12").UnMerge()
$ws.Range("C11:C12").UnMerge()
$ws.Range("E11:E12").UnMerge()
$ws.Range("F11:F12").UnMerge()

# ---------------------------------------------------------------------------
# 2. Row 11 (Eri Jesús Ocampo Alvarez) - only F.DOCENTE changes.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 0

# ---------------------------------------------------------------------------
# 3. Row 12 now becomes its own data row for "Ana Cortez Ramos" (previously
#    only D12 held a value as part of the merge).  Give the new cells the
#    same bordered / centred look as the rest of the data rows.
# ---------------------------------------------------------------------------
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = "Ana Cortez Ramos"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0

foreach ($addr in @("C12", "E12")) {
    $c = $ws.Range($addr)
    $c.Borders.LineStyle = 1
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
}

# ---------------------------------------------------------------------------
# 4. Row 13 now holds what used to be "Karla Flores Torres" (row 14) data.
# ---------------------------------------------------------------------------
$ws.Range("C13").Value = "Karla Flores Torres"
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 1

# ---------------------------------------------------------------------------
# 5. Old row 14 (now redundant) is cleared, and the TOTAL row moves from 15
#    up to 14, keeping the bold / bordered / centred TOTAL look.
# ---------------------------------------------------------------------------
$ws.Range("B14:C14").Clear()

$ws.Range("D14").Value = "TOTAL"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 2

foreach ($addr in @("D14", "E14", "F14")) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Calibri"
    $c.Font.Size = 11
    $c.Font.Bold = $true
    $c.Borders.LineStyle = 1
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
}

# Remove the now-empty old TOTAL row entirely (contents + formatting).
$ws.Range("D15:F15").Clear()

# ---------------------------------------------------------------------------
# 6. New signature block: name + job title, each merged across E:F, bold,
#    centred, no border.
# ---------------------------------------------------------------------------
$ws.Range("E18").Value = "MARITZA FLORES SARABIA"
$ws.Range("E19").Value = "JEFE DEL DEPARTAMENTO DE DESARROLLO ACADEMICO"

foreach ($addr in @("E18", "F18", "E19", "F19")) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Calibri"
    $c.Font.Size = 11
    $c.Font.Bold = $true
    $c.Borders.LineStyle = -4142
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
}

$ws.Range("E18:F18").Merge()
$ws.Range("E19:F19").Merge()

# ---------------------------------------------------------------------------
# 7. Column width tweaks.
# ---------------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 31.25
$ws.Columns("C").ColumnWidth = 23.43359375
$ws.Columns("E").ColumnWidth = 31.25
